# Update the "取得日時" (retrieved-at) timestamp column on the ランサーズ sheet.
# All rows that currently hold the old timestamp "2025-11-07 12:35:55"
# are refreshed to the new run's timestamp "2025-11-07 12:46:44".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-07 12:35:55"
$newTimestamp = "2025-11-07 12:46:44"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
